# Insert a new data row at row 215 (shifts existing rows 215-271 down to 216-272)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(215).Insert()

$ws.Cells.Item(215, 1).Value = 3
$ws.Cells.Item(215, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(215, 3).Value = "Coquimbo"
$ws.Cells.Item(215, 4).Value = 44551
$ws.Cells.Item(215, 5).Value = 5
$ws.Cells.Item(215, 6).Value = 100112040
$ws.Cells.Item(215, 7).Value = "Cilantro"
$ws.Cells.Item(215, 8).Value = "Sin especificar"
$ws.Cells.Item(215, 9).Value = "Primera"
$ws.Cells.Item(215, 10).Value = 140
$ws.Cells.Item(215, 11).Value = 5000
$ws.Cells.Item(215, 12).Value = 5500
$ws.Cells.Item(215, 13).Value = 5250
$ws.Cells.Item(215, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(215, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(215, 16).Value = 1750
$ws.Cells.Item(215, 17).Value = 3
$ws.Cells.Item(215, 18).Value = "Hortaliza"
